# Applies the cryptos-list refresh described in the commit diff.
# D-column values that look like plain numbers are written with a leading
# apostrophe so Excel stores them as text (matching the source workbook,
# where every Price cell is inline text, incl. ones like "1.00").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.558.36'
$ws.Range("E2").Value = '  +3.19%  '

$ws.Range("D3").Value = '1.841.38'

$ws.Range("D5").Value = "'231.79"
$ws.Range("E5").Value = '  +3.34%  '

$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = '  +3.10%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").Value = "'43.71"
$ws.Range("E8").Value = '  +11.05%  '

$ws.Range("E9").Value = '  +8.23%  '

$ws.Range("D10").Value = "'0.0704"
$ws.Range("E10").Value = '  +5.65%  '

$ws.Range("E11").Value = '  +2.58%  '

$ws.Range("D12").Value = '2.109.34'
$ws.Range("E12").Value = '  +2.18%  '

$ws.Range("D13").Value = '1.845.18'
$ws.Range("E13").Value = '  +2.31%  '

$ws.Range("D14").Value = "'11.31"
$ws.Range("E14").Value = '  +3.47%  '

$ws.Range("D15").Value = "'0.674"
$ws.Range("E15").Value = '  +7.07%  '

$ws.Range("D16").Value = "'4.73"
$ws.Range("E16").Value = '  +8.40%  '

$ws.Range("D17").Value = '35.535.60'
$ws.Range("E17").Value = '  +3.16%  '

$ws.Range("D18").Value = "'70.41"
$ws.Range("E18").Value = '  +3.57%  '

$ws.Range("D19").Value = '0.0₃0803'
$ws.Range("E19").Value = '  +4.83%  '

$ws.Range("D20").Value = "'245.03"
$ws.Range("E20").Value = '  +2.65%  '

$ws.Range("D21").Value = "'12.06"
$ws.Range("E21").Value = '  +8.57%  '

$ws.Range("D22").Value = "'4.62"
$ws.Range("E22").Value = '  +13.52%  '

$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("E24").Value = '  +2.80%  '

$ws.Range("D25").Value = "'171.79"
$ws.Range("E25").Value = '  +0.52%  '

$ws.Range("D26").Value = "'7.99"
$ws.Range("E26").Value = '  +4.58%  '

$ws.Range("E27").Value = '  +1.70%  '

$ws.Range("E28").Value = '  +0.72%  '

$ws.Range("E29").Value = '  +26.58%  '

$ws.Range("E30").Value = '  +0.17%  '

$ws.Range("D31").Value = '3.338.51'
$ws.Range("E31").Value = '  +37.41%  '

$ws.Range("E32").Value = '  +7.98%  '

$ws.Range("E33").Value = '  +7.34%  '

$ws.Range("E34").Value = '  +5.37%  '

$ws.Range("E35").Value = '  +1.65%  '

$ws.Range("D36").Value = "'95.49"
$ws.Range("E36").Value = '  +16.95%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = "'0.691"
$ws.Range("E37").Value = '  +8.08%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = "'1.14"
$ws.Range("E38").Value = '  +7.94%  '

$ws.Range("D39").Value = '1.346.94'
$ws.Range("E39").Value = '  +3.60%  '

$ws.Range("D40").Value = "'2.46"
$ws.Range("E40").Value = '  +7.54%  '

$ws.Range("D41").Value = "'15.45"
$ws.Range("E41").Value = '  +11.22%  '

$ws.Range("E42").Value = '  +5.31%  '

$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = '  +7.32%  '

$ws.Range("E44").Value = '  +4.48%  '

$ws.Range("E45").Value = '  +0.71%  '

$ws.Range("E46").Value = '  +0.93%  '

$ws.Range("E47").Value = '  +10.14%  '

$ws.Range("D48").Value = "'0.0519"
$ws.Range("E48").Value = '  +0.69%  '

$ws.Range("E49").Value = '  +2.45%  '

$ws.Range("E50").Value = '  +0.20%  '

$ws.Range("D51").Value = "'102.49"
$ws.Range("E51").Value = '  +0.71%  '
